$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 101
$ws.Range("D2").Value = 87
$ws.Range("E2").Value = 0.8613861386138614
$ws.Range("F2").Value = 0.8529411764705882
$ws.Range("G2").Value = 0.09847934560753153
$ws.Range("H2").Value = 0.08409512811622789
$ws.Range("I2").Value = 461526.6616369847
$ws.Range("J2").Value = 168149.3310194924
$ws.Range("L2").Value = 168149.3310194924
$ws.Range("M2").Value = 629675.9926564771
$ws.Range("N2").Value = 10098564.6888
$ws.Range("O2").Value = 9690823.758699998
$ws.Range("P2").Value = 0.01665081486342128
$ws.Range("Q2").Value = 0.01735139707483951
$ws.Range("G3").Value = 0.09664744730659888
$ws.Range("H3").Value = 0.08351090107075052
$ws.Range("I3").Value = 480378.4435696948
$ws.Range("J3").Value = 174976.8019918774
$ws.Range("L3").Value = 174976.8019918774
$ws.Range("M3").Value = 655355.2455615721
$ws.Range("N3").Value = 10549636.406764
$ws.Range("O3").Value = 10142263.248761
$ws.Range("P3").Value = 0.0165860504803454
$ws.Range("Q3").Value = 0.01725224416880058
$ws.Range("C4").Value = 104
$ws.Range("D4").Value = 90
$ws.Range("E4").Value = 0.8653846153846154
$ws.Range("F4").Value = 0.8653846153846154
$ws.Range("G4").Value = 0.09672223673719484
$ws.Range("H4").Value = 0.08370193563795707
$ws.Range("I4").Value = 511686.5069044705
$ws.Range("J4").Value = 183544.4610654761
$ws.Range("L4").Value = 183544.4610654761
$ws.Range("M4").Value = 695230.9679699468
$ws.Range("N4").Value = 10895668.70626692
$ws.Range("O4").Value = 10487324.35352383
$ws.Range("P4").Value = 0.01684563527155573
$ws.Range("Q4").Value = 0.01750155281540459
$ws.Range("G5").Value = 0.09467907571674584
$ws.Range("H5").Value = 0.08205519895451308
$ws.Range("I5").Value = 524337.9468682207
$ws.Range("J5").Value = 187296.2172757485
$ws.Range("L5").Value = 187296.2172757485
$ws.Range("M5").Value = 711634.1641439691
$ws.Range("N5").Value = 11320413.25095493
$ws.Range("O5").Value = 10909718.56762955
$ws.Range("P5").Value = 0.01654499823669858
$ws.Range("Q5").Value = 0.01716783215943617
$ws.Range("C6").Value = 105
$ws.Range("D6").Value = 91
$ws.Range("E6").Value = 0.8666666666666667
$ws.Range("F6").Value = 0.8584905660377359
$ws.Range("G6").Value = 0.09657166851730856
$ws.Range("H6").Value = 0.08300020599127433
$ws.Range("I6").Value = 550804.6883894347
$ws.Range("J6").Value = 197468.4053516046
$ws.Range("L6").Value = 197468.4053516046
$ws.Range("M6").Value = 748273.0937410393
$ws.Range("N6").Value = 11614789.51518358
$ws.Range("O6").Value = 11200323.99135843
$ws.Range("P6").Value = 0.01700146223859344
$ws.Range("Q6").Value = 0.01763059760628002
